$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date value that was updated for rows 2-17
# from 2023-10-25 (serial 45224) to 2023-11-03 (serial 45233).
$ws.Range("C2:C17").Value = 45233
